$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 130
$ws.Range("I11").Value = 130
$ws.Range("K11").Value = 130
$ws.Range("M11").Value = 10

$ws.Range("H96").Value = 1336
$ws.Range("I96").Value = 905
$ws.Range("K96").Value = 2715
$ws.Range("M96").Value = -1342

$ws.Range("H123").Value = 100780
$ws.Range("J123").Value = 100780
$ws.Range("L123").Value = 100780
$ws.Range("N123").Value = -110580

$ws.Range("H131").Value = 8343
$ws.Range("I131").Value = 3372
$ws.Range("K131").Value = 10116
$ws.Range("M131").Value = -5076

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 500
$ws.Range("I3").Value = 500
$ws.Range("K3").Value = 500
$ws.Range("M3").Value = -385

$ws.Range("H45").Value = 1941.6666
$ws.Range("I45").Value = 1912.5
$ws.Range("K45").Value = 1912.5
$ws.Range("M45").Value = -1535.5

$ws.Range("H61").Value = 1522.5454
$ws.Range("I61").Value = 1522.5454
$ws.Range("K61").Value = 1522.5454
$ws.Range("M61").Value = -1310.5454

$ws.Range("H74").Value = 5139.846
$ws.Range("I74").Value = 1970.6666
$ws.Range("J74").Value = 6090.6
$ws.Range("K74").Value = 1970.6666
$ws.Range("L74").Value = 6090.6
$ws.Range("M74").Value = -1096.6666
$ws.Range("N74").Value = -7838.6

$ws.Range("H77").Value = 5139.846
$ws.Range("I77").Value = 1970.6666
$ws.Range("J77").Value = 6090.6
$ws.Range("K77").Value = 9853.333000000001
$ws.Range("L77").Value = 30453
$ws.Range("M77").Value = -5485.333000000001
$ws.Range("N77").Value = -39189

$ws.Range("H122").Value = 387392.06
$ws.Range("I122").Value = 557517.7
$ws.Range("J122").Value = 4609.5
$ws.Range("K122").Value = 1672553.1
$ws.Range("L122").Value = 13828.5
$ws.Range("M122").Value = -1670103.1
$ws.Range("N122").Value = -18728.5

$ws.Range("H136").Value = 1522.5454
$ws.Range("I136").Value = 1522.5454
$ws.Range("K136").Value = 4567.6362
$ws.Range("M136").Value = -2017.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1633
$ws.Range("I5").Value = 1633
$ws.Range("K5").Value = 1633
$ws.Range("M5").Value = -1520

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4777.5884
$ws.Range("I31").Value = 2020
$ws.Range("J31").Value = 5626.077
$ws.Range("K31").Value = 2020
$ws.Range("L31").Value = 5626.077
$ws.Range("M31").Value = -1725
$ws.Range("N31").Value = -6216.077

$ws.Range("H34").Value = 4777.5884
$ws.Range("I34").Value = 2020
$ws.Range("J34").Value = 5626.077
$ws.Range("K34").Value = 2020
$ws.Range("L34").Value = 5626.077
$ws.Range("M34").Value = -1818
$ws.Range("N34").Value = -6030.077

$ws.Range("H69").Value = 17906.857
$ws.Range("I69").Value = 17906.857
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 17906.857
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -17157.857
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 17906.857
$ws.Range("I72").Value = 17906.857
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 53720.571
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -49976.571
$ws.Range("N72").ClearContents()

$ws.Range("H93").Value = 3816.3333
$ws.Range("I93").Value = 3816.3333
$ws.Range("K93").Value = 3816.3333
$ws.Range("M93").Value = -1944.3333

$ws.Range("H103").Value = 14736
$ws.Range("I103").Value = 12828
$ws.Range("K103").Value = 12828
$ws.Range("M103").Value = -11656

$ws.Range("H132").Value = 3836.125
$ws.Range("I132").Value = 3094.9412
$ws.Range("K132").Value = 9284.8236
$ws.Range("M132").Value = -6754.8236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2941841.8
$ws.Range("I4").Value = 2941841.8
$ws.Range("K4").Value = 8825525.399999999
$ws.Range("M4").Value = -8825413.399999999

$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 221
$ws.Range("N9").ClearContents()

$ws.Range("H12").Value = 28.285715
$ws.Range("J12").Value = 30.714285
$ws.Range("L12").Value = 92.142855
$ws.Range("N12").Value = -438.142855

$ws.Range("H32").Value = 14498.75
$ws.Range("I32").Value = 2998.6667
$ws.Range("J32").Value = 48999
$ws.Range("K32").Value = 8996.000100000001
$ws.Range("L32").Value = 146997
$ws.Range("M32").Value = -8713.000100000001
$ws.Range("N32").Value = -147563

$ws.Range("H37").Value = 109962.25
$ws.Range("J37").Value = 109962.25
$ws.Range("L37").Value = 329886.75
$ws.Range("N37").Value = -330110.75

$ws.Range("H118").Value = 1399.6666
$ws.Range("J118").Value = 1319.8
$ws.Range("L118").Value = 3959.4
$ws.Range("N118").Value = -6445.4

$ws.Range("H129").Value = 4749.75
$ws.Range("I129").Value = 6999
$ws.Range("K129").Value = 20997
$ws.Range("M129").Value = -15997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 18250
$ws.Range("J29").Value = 19000
$ws.Range("L29").Value = 19000
$ws.Range("N29").Value = -19580

$ws.Range("H107").Value = 1381.2273
$ws.Range("I107").Value = 2638.4
$ws.Range("K107").Value = 2638.4
$ws.Range("M107").Value = -718.4000000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7830.3335
$ws.Range("I61").Value = 9329
$ws.Range("J61").Value = 6331.6665
$ws.Range("K61").Value = 9329
$ws.Range("L61").Value = 6331.6665
$ws.Range("M61").Value = -9127
$ws.Range("N61").Value = -6735.6665

$ws.Range("H113").Value = 7830.3335
$ws.Range("I113").Value = 9329
$ws.Range("J113").Value = 6331.6665
$ws.Range("K113").Value = 9329
$ws.Range("L113").Value = 6331.6665
$ws.Range("M113").Value = -7159
$ws.Range("N113").Value = -10671.6665

$ws.Range("H132").Value = 5104.346
$ws.Range("I132").Value = 3588.2727
$ws.Range("J132").Value = 6216.1333
$ws.Range("K132").Value = 10764.8181
$ws.Range("L132").Value = 18648.3999
$ws.Range("M132").Value = -8234.8181
$ws.Range("N132").Value = -23708.3999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 25090.334
$ws.Range("J69").Value = 25090.334
$ws.Range("L69").Value = 25090.334
$ws.Range("N69").Value = -26588.334

$ws.Range("H72").Value = 25090.334
$ws.Range("J72").Value = 25090.334
$ws.Range("L72").Value = 75271.00199999999
$ws.Range("N72").Value = -82759.00199999999

$ws.Range("H96").Value = 1749
$ws.Range("I96").Value = 2500
$ws.Range("J96").Value = 998
$ws.Range("K96").Value = 2500
$ws.Range("L96").Value = 998
$ws.Range("M96").Value = -1127
$ws.Range("N96").Value = -3744

$ws.Range("H113").Value = 3134.3333
$ws.Range("I113").Value = 2250
$ws.Range("K113").Value = 6750
$ws.Range("M113").Value = -4580

$ws.Range("H132").Value = 1492
$ws.Range("I132").Value = 1031.2354
$ws.Range("K132").Value = 3093.7062
$ws.Range("M132").Value = -563.7062000000001

$ws.Range("H136").Value = 48499.547
$ws.Range("I136").Value = 2666.9412
$ws.Range("K136").Value = 8000.823600000001
$ws.Range("M136").Value = -5450.823600000001
